$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '88.078.53'
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").Value = '3.110.58'
$ws.Range("E3").Value = '  -2.15%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").Value = "'213.52"
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("D6").Value = "'633.05"
$ws.Range("E6").Value = '  +3.04%  '
$ws.Range("D7").Value = "'0.379"
$ws.Range("E7").Value = '  -2.76%  '
$ws.Range("D8").Value = "'0.811"
$ws.Range("E8").Value = '  +17.52%  '
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("D10").Value = '3.106.32'
$ws.Range("E10").Value = '  -2.10%  '
$ws.Range("D11").Value = "'0.589"
$ws.Range("E11").Value = '  +2.26%  '
$ws.Range("D12").Value = "'0.179"
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").Value = "'0.0000244"
$ws.Range("E13").Value = '  -3.86%  '
$ws.Range("D14").Value = "'5.34"
$ws.Range("E14").Value = '  +1.90%  '
$ws.Range("D15").Value = '87.812.41'
$ws.Range("E15").Value = '  -1.95%  '
$ws.Range("D16").Value = '3.675.16'
$ws.Range("E16").Value = '  -2.51%  '
$ws.Range("D17").Value = "'31.95"
$ws.Range("E17").Value = '  -2.64%  '
$ws.Range("D18").Value = '3.103.08'
$ws.Range("E18").Value = '  -3.34%  '
$ws.Range("D19").Value = "'3.44"
$ws.Range("E19").Value = '  +4.95%  '
$ws.Range("D20").Value = "'0.0000215"
$ws.Range("E20").Value = '  +9.71%  '
$ws.Range("D21").Value = "'13.28"
$ws.Range("E21").Value = '  -0.69%  '
$ws.Range("D22").Value = "'423.07"
$ws.Range("E22").Value = '  -2.93%  '
$ws.Range("D23").Value = "'8.39"
$ws.Range("E23").Value = '  -2.16%  '
$ws.Range("D24").Value = "'4.86"
$ws.Range("E24").Value = '  -3.86%  '
$ws.Range("D25").Value = "'5.49"
$ws.Range("E25").Value = '  +7.03%  '
$ws.Range("D26").Value = "'83.53"
$ws.Range("E26").Value = '  +10.97%  '
$ws.Range("D27").Value = "'11.40"
$ws.Range("E27").Value = '  -2.50%  '
$ws.Range("D28").Value = '3.269.43'
$ws.Range("E28").Value = '  -2.35%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").Value = "'0.152"
$ws.Range("E31").Value = '  -8.78%  '
$ws.Range("D32").Value = "'8.15"
$ws.Range("E32").Value = '  -3.14%  '
$ws.Range("D33").Value = "'3.79"
$ws.Range("E33").Value = '  -5.84%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").Value = "'503.55"
$ws.Range("E34").Value = '  -5.55%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = "'0.146"
$ws.Range("E35").Value = '  +14.42%  '
$ws.Range("D36").Value = "'6.80"
$ws.Range("E36").Value = '  -3.01%  '
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("E38").Value = '  -2.10%  '
$ws.Range("D39").Value = "'22.43"
$ws.Range("E39").Value = '  +2.07%  '
$ws.Range("D40").Value = "'22.19"
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").Value = "'0.366"
$ws.Range("E43").Value = '  -1.47%  '
$ws.Range("D44").Value = "'1.85"
$ws.Range("E44").Value = '  -3.74%  '
$ws.Range("E45").Value = '  +10.07%  '
$ws.Range("D46").Value = "'146.11"
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("D47").Value = "'43.63"
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").Value = "'0.0661"
$ws.Range("E48").Value = '  +12.04%  '
$ws.Range("D49").Value = "'161.17"
$ws.Range("E49").Value = '  -6.43%  '
$ws.Range("D50").Value = "'0.716"
$ws.Range("E50").Value = '  +2.22%  '
$ws.Range("D51").Value = "'1.19"
$ws.Range("E51").Value = '  -3.29%  '
